# Slide "ERD":
#  - remove the now-empty "body" placeholder shape (it only ever held the
#    placeholder text "그림")
#  - reposition/resize the ERD diagram picture

$p = $ppt.ActivePresentation

# --- locate the "ERD" slide (defaults to slide 14, where it lives in this
# deck, but searches by title text too so the script keeps working if the
# slide order ever changes) ---
$s = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        $isTitle = $false
        try { $isTitle = $shp.PlaceholderFormat.Type -eq 1 } catch { $isTitle = $false }
        if ($isTitle -and $shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "ERD") {
            $s = $slide
        }
    }
}
if ($s -eq $null) {
    $s = $p.Slides.Item(14)
}

# --- find the index of the empty body placeholder ("그림") and of the picture ---
$bodyIdx = -1
$picIdx = -1
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    $isPlaceholder = $false
    try { $isPlaceholder = $shp.Type -eq 14 } catch { $isPlaceholder = $false }
    if ($isPlaceholder -and $shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "그림") {
        $bodyIdx = $i
    }
    if (-not $shp.HasTextFrame) {
        $picIdx = $i
    }
}

# A placeholder shape needs two Delete() calls to be fully removed from the
# slide: the first clears its contents (PowerPoint resets it to a fresh,
# empty generic placeholder), the second removes the now-empty shape from
# the shape tree. Re-fetch the shape from the collection each time instead
# of re-using a cached reference, since the first Delete() invalidates it.
if ($bodyIdx -gt 0) {
    $s.Shapes.Item($bodyIdx).Delete()
    $s.Shapes.Item($bodyIdx).Delete()
    if ($picIdx -gt $bodyIdx) {
        $picIdx = $picIdx - 1
    }
}

# --- reposition / resize the ERD picture ---
# Shape.Left/Top/Width/Height are expressed in points; the target geometry
# below is specified in EMU (1 pt = 12700 EMU = 1/72 in = 1/914400 in-EMU).
if ($picIdx -gt 0) {
    $pic = $s.Shapes.Item($picIdx)
    $pic.Left = 2219417 / 914400 * 72
    $pic.Top = 146481 / 914400 * 72
    $pic.Width = 8966447 / 914400 * 72
    $pic.Height = 6565037 / 914400 * 72
}

Write-Output "edit complete"
